$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1431
$ws.Range("F3").Value = 1408
$ws.Range("F4").Value = 404
$ws.Range("F5").Value = 221
$ws.Range("F6").Value = 667
$ws.Range("F7").Value = 25
$ws.Range("F8").Value = 613
$ws.Range("F10").Value = 73
$ws.Range("F11").Value = 1362
$ws.Range("F12").Value = 31742
$ws.Range("F13").Value = 6840
$ws.Range("F14").Value = 103
$ws.Range("F15").Value = 341
$ws.Range("F16").Value = 562
$ws.Range("F17").Value = 424
$ws.Range("F21").Value = 427
$ws.Range("F24").Value = 313
$ws.Range("F25").Value = 374
$ws.Range("F26").Value = 421
$ws.Range("F28").Value = 181
$ws.Range("F29").Value = 42
$ws.Range("F30").Value = 724
$ws.Range("F31").Value = 281
$ws.Range("F33").Value = 711
$ws.Range("F34").Value = 103
$ws.Range("F35").Value = 41
$ws.Range("F36").Value = 775
$ws.Range("F37").Value = 278
$ws.Range("F38").Value = 48
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 1140
$ws.Range("F5").Value = 144
$ws.Range("F7").Value = 4314
$ws.Range("F19").Value = 4281
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1417
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1417
$ws.Range("F4").Value = 1140
$ws.Range("F5").Value = 1431
$ws.Range("F6").Value = 1408
$ws.Range("F7").Value = 221
$ws.Range("F8").Value = 667
$ws.Range("F9").Value = 25
$ws.Range("F10").Value = 613
$ws.Range("F12").Value = 73
$ws.Range("F13").Value = 1362
$ws.Range("F14").Value = 144
$ws.Range("F20").Value = 6840
$ws.Range("F21").Value = 103
$ws.Range("F22").Value = 341
$ws.Range("F24").Value = 562
$ws.Range("F25").Value = 424
$ws.Range("F31").Value = 427
$ws.Range("F34").Value = 313
$ws.Range("F35").Value = 374
$ws.Range("F36").Value = 421
$ws.Range("F38").Value = 181
$ws.Range("F39").Value = 42
$ws.Range("F40").Value = 724
$ws.Range("F42").Value = 281
$ws.Range("F44").Value = 103
$ws.Range("F45").Value = 775
$ws.Range("F46").Value = 278
$ws.Range("F47").Value = 48
